$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 129; this shifts existing rows 129-160 down to 130-161
$ws.Rows.Item(129).Insert()

# Populate the newly inserted row 129 with a fresh weekly price observation.
# (Same market/category metadata as its neighboring rows; new date & price figures.)
$ws.Cells.Item(129, 1).Value = 11
$ws.Cells.Item(129, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(129, 3).Value = "Bíobío"
$ws.Cells.Item(129, 4).Value = 44663
$ws.Cells.Item(129, 5).Value = 8
$ws.Cells.Item(129, 6).Value = 100112003
$ws.Cells.Item(129, 7).Value = "Ajo"
$ws.Cells.Item(129, 8).Value = "Chino"
$ws.Cells.Item(129, 9).Value = "Primera"
$ws.Cells.Item(129, 10).Value = 400
$ws.Cells.Item(129, 11).Value = 22000
$ws.Cells.Item(129, 12).Value = 24000
$ws.Cells.Item(129, 13).Value = 23000
$ws.Cells.Item(129, 14).Value = "$/caja 10 kilos"
$ws.Cells.Item(129, 15).Value = "China"
$ws.Cells.Item(129, 16).Value = 2300
$ws.Cells.Item(129, 17).Value = 10
$ws.Cells.Item(129, 18).Value = "Hortaliza"
